$wb = $excel.ActiveWorkbook

# ALC row 17: One for the Road
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1433525.9
$ws.Range("J17").Value = 1433525.9
$ws.Range("L17").Value = 4300577.699999999
$ws.Range("N17").Value = -4300913.699999999

# ALC row 39: Riches' Brew
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 227
$ws.Range("I39").Value = 226.89473
$ws.Range("J39").Value = 227.5
$ws.Range("K39").Value = 680.6841900000001
$ws.Range("L39").Value = 682.5
$ws.Range("M39").Value = -384.6841900000001
$ws.Range("N39").Value = -1274.5

# ALC row 40: Stuck in the Moment
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 11995
$ws.Range("J40").Value = 11995
$ws.Range("L40").Value = 11995
$ws.Range("N40").Value = -12345

# ALC row 55: A Real Smooth Move
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 266.15384
$ws.Range("I55").Value = 176.5
$ws.Range("K55").Value = 176.5
$ws.Range("M55").Value = 37.5

# ALC row 62: The Mustache Suits Him
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8122.6
$ws.Range("J62").Value = 12860.429
$ws.Range("L62").Value = 12860.429
$ws.Range("N62").Value = -14108.429

# ALC row 65: Forgery of Convenience (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 8122.6
$ws.Range("J65").Value = 12860.429
$ws.Range("L65").Value = 64302.145
$ws.Range("N65").Value = -70542.145

# ALC row 106: Making Your Mark
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 8545
$ws.Range("I106").Value = 4170.857
$ws.Range("K106").Value = 4170.857
$ws.Range("M106").Value = -3539.857

# ALC row 132: Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 31718.4
$ws.Range("J132").Value = 36668.668
$ws.Range("L132").Value = 110006.004
$ws.Range("N132").Value = -115066.004

# ALC row 140: Tome for Tradition
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 89194.5
$ws.Range("J140").Value = 89194.5
$ws.Range("L140").Value = 89194.5
$ws.Range("N140").Value = -99554.5

# ARM row 2: Ain't Got No Ingots
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4186.3335
$ws.Range("I2").Value = 973.3333
$ws.Range("J2").Value = 7399.3335
$ws.Range("K2").Value = 973.3333
$ws.Range("L2").Value = 7399.3335
$ws.Range("M2").Value = -860.3333
$ws.Range("N2").Value = -7625.3335

# ARM row 33: A Leg to Stand On
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").ClearContents()
$ws.Range("N33").Value = 0

# ARM row 45: Hollow Hallmarks
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1557.3529
$ws.Range("I45").Value = 1422.25
$ws.Range("J45").Value = 1881.6
$ws.Range("K45").Value = 1422.25
$ws.Range("L45").Value = 1881.6
$ws.Range("M45").Value = -1045.25
$ws.Range("N45").Value = -2635.6

# ARM row 61: Dealing with the Tough Stuff
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17866.541
$ws.Range("I61").Value = 3682
$ws.Range("K61").Value = 3682
$ws.Range("M61").Value = -3470

# ARM row 116: No Scope
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 4186.3335
$ws.Range("I116").Value = 973.3333
$ws.Range("J116").Value = 7399.3335
$ws.Range("K116").Value = 973.3333
$ws.Range("L116").Value = 7399.3335
$ws.Range("M116").Value = 1320.6667
$ws.Range("N116").Value = -11987.3335

# ARM row 122: Haste for High Durium
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4623.923
$ws.Range("I122").Value = 1674.5
$ws.Range("J122").Value = 5934.778
$ws.Range("K122").Value = 5023.5
$ws.Range("L122").Value = 17804.334
$ws.Range("M122").Value = -2573.5
$ws.Range("N122").Value = -22704.334

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3587508.5
$ws.Range("I132").Value = 5898.8887
$ws.Range("K132").Value = 17696.6661
$ws.Range("M132").Value = -15166.6661

# ARM row 136: Metal with Mettle
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 17866.541
$ws.Range("I136").Value = 3682
$ws.Range("K136").Value = 11046
$ws.Range("M136").Value = -8496

# BSM row 3: Hells Bells
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4186.3335
$ws.Range("I3").Value = 973.3333
$ws.Range("J3").Value = 7399.3335
$ws.Range("K3").Value = 973.3333
$ws.Range("L3").Value = 7399.3335
$ws.Range("M3").Value = -859.3333
$ws.Range("N3").Value = -7627.3335

# BSM row 54: Get Me to the War on Time
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 591.25
$ws.Range("I54").Value = 591.25
$ws.Range("K54").Value = 591.25
$ws.Range("M54").Value = -107.25

# BSM row 134: Ruthenium Supremium
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 10742.692
$ws.Range("I134").Value = 1936.625
$ws.Range("K134").Value = 5809.875
$ws.Range("M134").Value = -3274.875

# CRP row 22: Driving Up the Wall
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3249.0833
$ws.Range("I22").Value = 2874.25
$ws.Range("J22").Value = 3436.5
$ws.Range("K22").Value = 2874.25
$ws.Range("L22").Value = 3436.5
$ws.Range("M22").Value = -2524.25
$ws.Range("N22").Value = -4136.5

# CRP row 31: Wall Not Found
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29002.736
$ws.Range("I31").Value = 19241.334
$ws.Range("J31").Value = 33508
$ws.Range("K31").Value = 19241.334
$ws.Range("L31").Value = 33508
$ws.Range("M31").Value = -18946.334
$ws.Range("N31").Value = -34098

# CRP row 34: Armoires of the Rich and Famous
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 29002.736
$ws.Range("I34").Value = 19241.334
$ws.Range("J34").Value = 33508
$ws.Range("K34").Value = 19241.334
$ws.Range("L34").Value = 33508
$ws.Range("M34").Value = -19039.334
$ws.Range("N34").Value = -33912

# CRP row 58: You Do the Heavy Lifting
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 15148.719
$ws.Range("I58").Value = 6602
$ws.Range("K58").Value = 6602
$ws.Range("M58").Value = -6399

# CRP row 105: Zelkova, My Love
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 22843.857
$ws.Range("I105").Value = 50954.5
$ws.Range("J105").Value = 11599.6
$ws.Range("K105").Value = 50954.5
$ws.Range("L105").Value = 11599.6
$ws.Range("M105").Value = -49207.5
$ws.Range("N105").Value = -15093.6

# CRP row 107: Built to Last
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2888.6897
$ws.Range("I107").Value = 1289.7368
$ws.Range("J107").Value = 5926.7
$ws.Range("K107").Value = 1289.7368
$ws.Range("L107").Value = 5926.7
$ws.Range("M107").Value = 630.2632000000001
$ws.Range("N107").Value = -9766.700000000001

# CRP row 132: Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 10642.25
$ws.Range("I132").Value = 2973.5557
$ws.Range("J132").Value = 20502
$ws.Range("K132").Value = 8920.667099999999
$ws.Range("L132").Value = 61506
$ws.Range("M132").Value = -6390.667099999999
$ws.Range("N132").Value = -66566

# CRP row 134: Wood You Be Quiet
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 50010160
$ws.Range("I134").Value = 2384.8333
$ws.Range("J134").Value = 71442060
$ws.Range("K134").Value = 7154.499899999999
$ws.Range("L134").Value = 214326180
$ws.Range("M134").Value = -4619.499899999999
$ws.Range("N134").Value = -214331250

# CRP row 136: Turali Quality
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 15148.719
$ws.Range("I136").Value = 6602
$ws.Range("K136").Value = 19806
$ws.Range("M136").Value = -17256

# CUL row 70: Persona non Gratin
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 24999.5
$ws.Range("I70").Value = 24999
$ws.Range("K70").Value = 74997
$ws.Range("M70").Value = -74682

# CUL row 73: Recipe for Disaster (L)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 24999.5
$ws.Range("I73").Value = 24999
$ws.Range("K73").Value = 74997
$ws.Range("M73").Value = -73905

# CUL row 107: Slippery Service
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3914056
$ws.Range("I107").Value = 966.3333
$ws.Range("K107").Value = 2898.9999
$ws.Range("M107").Value = -978.9998999999998

# CUL row 113: Can't Eat Just One
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 8688.352999999999
$ws.Range("I113").Value = 18152.572
$ws.Range("K113").Value = 54457.716
$ws.Range("M113").Value = -52287.716

# CUL row 124: Bobbing for Compliments
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 2924.8572
$ws.Range("I124").Value = 2412.3333
$ws.Range("J124").Value = 6000
$ws.Range("K124").Value = 7236.999899999999
$ws.Range("L124").Value = 18000
$ws.Range("M124").Value = -2326.999899999999
$ws.Range("N124").Value = -27820

# GSM row 102: Put the Metal to the Peddle
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7552
$ws.Range("I102").Value = 5404.1816
$ws.Range("K102").Value = 5404.1816
$ws.Range("M102").Value = -3782.1816

# GSM row 122: Awarding Academic Excellence
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7740.615
$ws.Range("I122").Value = 1420.8182
$ws.Range("K122").Value = 4262.4546
$ws.Range("M122").Value = -1812.4546

# GSM row 123: Workplace Workout
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 91925
$ws.Range("J123").Value = 91925
$ws.Range("L123").Value = 91925
$ws.Range("N123").Value = -96825

# GSM row 135: Fan of the Foreign
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 153306.84
$ws.Range("J135").Value = 153306.84
$ws.Range("L135").Value = 153306.84
$ws.Range("N135").Value = -163446.84

# LTW row 7: Tan Before the Ban
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8356.632
$ws.Range("I7").Value = 7359.5415
$ws.Range("K7").Value = 7359.5415
$ws.Range("M7").Value = -7247.5415

# LTW row 58: Handle with Care
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 19500
$ws.Range("I58").Value = 19000
$ws.Range("J58").Value = 20000
$ws.Range("K58").Value = 19000
$ws.Range("L58").Value = 20000
$ws.Range("M58").Value = -18740
$ws.Range("N58").Value = -20520

# LTW row 122: Hell on Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7285.3228
$ws.Range("J122").Value = 10038.077
$ws.Range("L122").Value = 30114.231
$ws.Range("N122").Value = -35014.231

# LTW row 126: Battered Books
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 8356.632
$ws.Range("I126").Value = 7359.5415
$ws.Range("K126").Value = 22078.6245
$ws.Range("M126").Value = -19608.6245

# LTW row 132: Tenets of Tanning
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3100493
$ws.Range("I132").Value = 12400.167
$ws.Range("K132").Value = 37200.501
$ws.Range("M132").Value = -34670.501

# LTW row 136: Respect for Br'aax
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 11302.0205
$ws.Range("I136").Value = 10234.6
$ws.Range("K136").Value = 30703.8
$ws.Range("M136").Value = -28153.8

# WVR row 92: Modest Beginnings
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 39999
$ws.Range("J92").Value = 39999
$ws.Range("L92").Value = 39999
$ws.Range("N92").Value = -44991

# WVR row 132: Comfy Cabins
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8095.3335
$ws.Range("I132").Value = 3227.4092
$ws.Range("J132").Value = 14395
$ws.Range("K132").Value = 9682.2276
$ws.Range("L132").Value = 43185
$ws.Range("M132").Value = -7152.2276
$ws.Range("N132").Value = -48245
